# Scheduled-runner refresh of Anima_Profits data: updates cached
# market-board price/profit figures (columns H-N) on affected rows
# across the ALC/ARM/CRP/CUL/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 12680.833
$ws.Cells.Item(28, 9).Value = 216.4
$ws.Cells.Item(28, 11).Value = 216.4
$ws.Cells.Item(28, 13).Value = 268.6

$ws.Cells.Item(98, 8).Value = 1248
$ws.Cells.Item(98, 9).Value = 934.25
$ws.Cells.Item(98, 10).Value = 2503
$ws.Cells.Item(98, 11).Value = 934.25
$ws.Cells.Item(98, 12).Value = 2503
$ws.Cells.Item(98, 13).Value = 563.75
$ws.Cells.Item(98, 14).Value = -5499

$ws.Cells.Item(112, 8).Value = 5896.9814
$ws.Cells.Item(112, 10).Value = 6396.6733
$ws.Cells.Item(112, 12).Value = 19190.0199
$ws.Cells.Item(112, 14).Value = -21406.0199

$ws.Cells.Item(122, 8).Value = 1248
$ws.Cells.Item(122, 9).Value = 934.25
$ws.Cells.Item(122, 10).Value = 2503
$ws.Cells.Item(122, 11).Value = 2802.75
$ws.Cells.Item(122, 12).Value = 7509
$ws.Cells.Item(122, 13).Value = -352.75
$ws.Cells.Item(122, 14).Value = -12409

$ws.Cells.Item(138, 8).Value = 1776.44
$ws.Cells.Item(138, 9).Value = 666.06665
$ws.Cells.Item(138, 10).Value = 1972.3882
$ws.Cells.Item(138, 11).Value = 1998.19995
$ws.Cells.Item(138, 12).Value = 5917.1646
$ws.Cells.Item(138, 13).Value = 3141.80005
$ws.Cells.Item(138, 14).Value = -16197.1646

$ws.Cells.Item(141, 8).Value = 7323.5713
$ws.Cells.Item(141, 9).Value = 2998.3333
$ws.Cells.Item(141, 10).Value = 10567.5
$ws.Cells.Item(141, 11).Value = 8994.999899999999
$ws.Cells.Item(141, 12).Value = 31702.5
$ws.Cells.Item(141, 13).Value = -3814.999899999999
$ws.Cells.Item(141, 14).Value = -42062.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1868.4762
$ws.Cells.Item(74, 9).Value = 1749.3684
$ws.Cells.Item(74, 10).Value = 3000
$ws.Cells.Item(74, 11).Value = 1749.3684
$ws.Cells.Item(74, 12).Value = 3000
$ws.Cells.Item(74, 13).Value = -875.3684000000001
$ws.Cells.Item(74, 14).Value = -4748

$ws.Cells.Item(77, 8).Value = 1868.4762
$ws.Cells.Item(77, 9).Value = 1749.3684
$ws.Cells.Item(77, 10).Value = 3000
$ws.Cells.Item(77, 11).Value = 8746.842000000001
$ws.Cells.Item(77, 12).Value = 15000
$ws.Cells.Item(77, 13).Value = -4378.842000000001
$ws.Cells.Item(77, 14).Value = -23736

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 26916.834
$ws.Cells.Item(4, 10).Value = 26916.834
$ws.Cells.Item(4, 12).Value = 26916.834
$ws.Cells.Item(4, 14).Value = -27140.834

$ws.Cells.Item(31, 8).Value = 6930.7026
$ws.Cells.Item(31, 9).Value = 1536.88
$ws.Cells.Item(31, 10).Value = 18167.834
$ws.Cells.Item(31, 11).Value = 1536.88
$ws.Cells.Item(31, 12).Value = 18167.834
$ws.Cells.Item(31, 13).Value = -1241.88
$ws.Cells.Item(31, 14).Value = -18757.834

$ws.Cells.Item(34, 8).Value = 6930.7026
$ws.Cells.Item(34, 9).Value = 1536.88
$ws.Cells.Item(34, 10).Value = 18167.834
$ws.Cells.Item(34, 11).Value = 1536.88
$ws.Cells.Item(34, 12).Value = 18167.834
$ws.Cells.Item(34, 13).Value = -1334.88
$ws.Cells.Item(34, 14).Value = -18571.834

$ws.Cells.Item(86, 8).Value = 2706.88
$ws.Cells.Item(86, 9).Value = 2785.875
$ws.Cells.Item(86, 10).Value = 2566.4443
$ws.Cells.Item(86, 11).Value = 2785.875
$ws.Cells.Item(86, 12).Value = 2566.4443
$ws.Cells.Item(86, 13).Value = -1662.875
$ws.Cells.Item(86, 14).Value = -4812.4443

$ws.Cells.Item(89, 8).Value = 2706.88
$ws.Cells.Item(89, 9).Value = 2785.875
$ws.Cells.Item(89, 10).Value = 2566.4443
$ws.Cells.Item(89, 11).Value = 13929.375
$ws.Cells.Item(89, 12).Value = 12832.2215
$ws.Cells.Item(89, 13).Value = -8313.375
$ws.Cells.Item(89, 14).Value = -24064.2215

$ws.Cells.Item(132, 8).Value = 4763603.5
$ws.Cells.Item(132, 9).Value = 1267.6538
$ws.Cells.Item(132, 10).Value = 18521464
$ws.Cells.Item(132, 11).Value = 3802.9614
$ws.Cells.Item(132, 12).Value = 55564392
$ws.Cells.Item(132, 13).Value = -1272.9614
$ws.Cells.Item(132, 14).Value = -55569452

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 627.3889
$ws.Cells.Item(5, 9).Value = 580.8125
$ws.Cells.Item(5, 11).Value = 1742.4375
$ws.Cells.Item(5, 13).Value = -1630.4375

$ws.Cells.Item(35, 8).Value = 3881.2
$ws.Cells.Item(35, 10).Value = 4776.5
$ws.Cells.Item(35, 12).Value = 14329.5
$ws.Cells.Item(35, 14).Value = -14905.5

$ws.Cells.Item(113, 8).Value = 869.29034
$ws.Cells.Item(113, 9).Value = 472.30768
$ws.Cells.Item(113, 10).Value = 1156
$ws.Cells.Item(113, 11).Value = 1416.92304
$ws.Cells.Item(113, 12).Value = 3468
$ws.Cells.Item(113, 13).Value = 753.0769599999999
$ws.Cells.Item(113, 14).Value = -7808

$ws.Cells.Item(122, 8).Value = 17012.334
$ws.Cells.Item(122, 10).Value = 50149.5
$ws.Cells.Item(122, 12).Value = 451345.5
$ws.Cells.Item(122, 14).Value = -456245.5

$ws.Cells.Item(132, 8).Value = 2618.6365
$ws.Cells.Item(132, 10).Value = 2422.1428
$ws.Cells.Item(132, 12).Value = 21799.2852
$ws.Cells.Item(132, 14).Value = -26859.2852

$ws.Cells.Item(135, 8).Value = 627.3889
$ws.Cells.Item(135, 9).Value = 580.8125
$ws.Cells.Item(135, 11).Value = 5227.3125
$ws.Cells.Item(135, 13).Value = -2692.3125

$ws.Cells.Item(136, 8).Value = 3831.6667
$ws.Cells.Item(136, 9).Value = 2990
$ws.Cells.Item(136, 11).Value = 8970
$ws.Cells.Item(136, 13).Value = -3870

$ws.Cells.Item(137, 8).Value = 30983
$ws.Cells.Item(137, 9).Value = 57499.5
$ws.Cells.Item(137, 10).Value = 4466.5
$ws.Cells.Item(137, 11).Value = 172498.5
$ws.Cells.Item(137, 12).Value = 13399.5
$ws.Cells.Item(137, 13).Value = -167398.5
$ws.Cells.Item(137, 14).Value = -23599.5

$ws.Cells.Item(139, 8).Value = 3270.6365
$ws.Cells.Item(139, 9).Value = 4006
$ws.Cells.Item(139, 11).Value = 12018
$ws.Cells.Item(139, 13).Value = -6878

$ws.Cells.Item(140, 8).Value = 1585.7241
$ws.Cells.Item(140, 10).Value = 2021
$ws.Cells.Item(140, 12).Value = 6063
$ws.Cells.Item(140, 14).Value = -16423

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3348.4707
$ws.Cells.Item(68, 9).Value = 3133.3333
$ws.Cells.Item(68, 10).Value = 3394.5715
$ws.Cells.Item(68, 11).Value = 3133.3333
$ws.Cells.Item(68, 12).Value = 3394.5715
$ws.Cells.Item(68, 13).Value = -2384.3333
$ws.Cells.Item(68, 14).Value = -4892.5715

$ws.Cells.Item(71, 8).Value = 3348.4707
$ws.Cells.Item(71, 9).Value = 3133.3333
$ws.Cells.Item(71, 10).Value = 3394.5715
$ws.Cells.Item(71, 11).Value = 15666.6665
$ws.Cells.Item(71, 12).Value = 16972.8575
$ws.Cells.Item(71, 13).Value = -11922.6665
$ws.Cells.Item(71, 14).Value = -24460.8575

$ws.Cells.Item(82, 8).Value = 83336150
$ws.Cells.Item(82, 9).Value = 166668670
$ws.Cells.Item(82, 10).Value = 3633.3333
$ws.Cells.Item(82, 11).Value = 166668670
$ws.Cells.Item(82, 12).Value = 3633.3333
$ws.Cells.Item(82, 13).Value = -166668309
$ws.Cells.Item(82, 14).Value = -4355.3333

$ws.Cells.Item(85, 8).Value = 83336150
$ws.Cells.Item(85, 9).Value = 166668670
$ws.Cells.Item(85, 10).Value = 3633.3333
$ws.Cells.Item(85, 11).Value = 166668670
$ws.Cells.Item(85, 12).Value = 3633.3333
$ws.Cells.Item(85, 13).Value = -166667422
$ws.Cells.Item(85, 14).Value = -6129.3333

$ws.Cells.Item(122, 8).Value = 2309.0908
$ws.Cells.Item(122, 9).Value = 1566.6666
$ws.Cells.Item(122, 10).Value = 3200
$ws.Cells.Item(122, 11).Value = 4699.9998
$ws.Cells.Item(122, 12).Value = 9600
$ws.Cells.Item(122, 13).Value = -2249.9998
$ws.Cells.Item(122, 14).Value = -14500

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4168765.5
$ws.Cells.Item(132, 9).Value = 2022.3478
$ws.Cells.Item(132, 10).Value = 9806124
$ws.Cells.Item(132, 11).Value = 6067.0434
$ws.Cells.Item(132, 12).Value = 29418372
$ws.Cells.Item(132, 13).Value = -3537.0434
$ws.Cells.Item(132, 14).Value = -29423432
